# Refresh the cryptocurrency Price (D) / 1h Volume % (E) columns with the
# latest scrape values, matching the automated GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values (e.g. "1.001", "300.38") look like plain
# numbers, which Excel would silently convert from text to a Number on entry.
# The sheet stores every Price/Volume cell as text, so force those specific
# cells to a text format before writing the new value, preserving their type.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.466.24"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.646.54"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "300.38"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("D7").Value = "0.3791"
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("D8").Value = "50.56"
$ws.Range("E8").Value = "  -0.83%  "
$ws.Range("D9").Value = "0.3499"
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").Value = "1.216"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "22.08"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "6.306"
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").Value = "7.250"
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("D16").Value = "0.00001211"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "1.647.19"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "95.06"
$ws.Range("E18").Value = "  -2.48%  "
$ws.Range("D19").Value = "0.06966"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "6.618"
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  -1.54%  "
$ws.Range("D24").Value = "23.471.33"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "2.419"
$ws.Range("E25").Value = "  -3.26%  "
$ws.Range("D26").Value = "2.965"
$ws.Range("E26").Value = "  -1.83%  "
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").Value = "151.76"
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("D29").Value = "5.184"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("D30").Value = "131.71"
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").Value = "1.828.44"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").Value = "6.830"
$ws.Range("E32").Value = "  -3.90%  "
$ws.Range("D33").Value = "2.137"
$ws.Range("E33").Value = "  -4.72%  "
$ws.Range("D34").Value = "11.18"
$ws.Range("E34").Value = "  -7.72%  "
$ws.Range("D35").Value = "0.9866"
$ws.Range("E35").Value = "  -6.90%  "
$ws.Range("D36").Value = "0.02687"
$ws.Range("E36").Value = "  -3.79%  "
$ws.Range("D37").Value = "0.08779"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "5.911"
$ws.Range("E38").Value = "  -2.96%  "
$ws.Range("E40").Value = "  -2.94%  "
$ws.Range("D41").Value = "12.83"
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("D42").Value = "0.6868"
$ws.Range("E42").Value = "  -1.59%  "
$ws.Range("D43").Value = "1.294"
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("D44").Value = "15.64"
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").Value = "0.9998"
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D46").Value = "0.6378"
$ws.Range("E46").Value = "  -2.03%  "
$ws.Range("D47").Value = "3.928"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").Value = "2.244"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("D49").Value = "0.07678"
$ws.Range("E49").Value = "  -2.44%  "
$ws.Range("D50").Value = "127.07"
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("D51").Value = "1.230"
$ws.Range("E51").Value = "  +2.20%  "
